# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet (fund holdings detail) positioned
#    right before the "总计" (total) summary sheet - built by cloning the
#    existing "2021-Q4" sheet (same header/layout) and overwriting its data
#    rows.
# 2) Prepend a "2022-Q1" row to the "总计" summary sheet, shifting the
#    existing history rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet from the "2021-Q4" template
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $template)
$newSheet = $wb.Worksheets.Item($template.Index + 1)
$newSheet.Name = "2022-Q1"

# Columns B and D:G hold numeric-looking text ("002423", "3.62", ...) that
# must stay text, not get auto-coerced into numbers. Force text format
# before assigning, then strip the format change back off afterwards
# (paste-format-only from an already-unstyled cell) so the saved cell has
# no stray style index, matching the sibling quarter sheets exactly.
$newSheet.Range("B2:B4").NumberFormat = "@"
$newSheet.Range("D2:G4").NumberFormat = "@"

$newSheet.Range("B2").Value = "002423"
$newSheet.Range("C2").Value = "华宝兴业标普美国消费(QDII-LOF)美元"
$newSheet.Range("D2").Value = "3.62"
$newSheet.Range("E2").Value = "94.37"
$newSheet.Range("F2").Value = "2.32"
$newSheet.Range("G2").Value = "0.0840"
$newSheet.Range("H2").Value = 8

$newSheet.Range("B3").Value = "162415"
$newSheet.Range("C3").Value = "华宝标普美国消费(QDII-LOF)人民币A"
$newSheet.Range("D3").Value = "3.62"
$newSheet.Range("E3").Value = "94.37"
$newSheet.Range("F3").Value = "2.32"
$newSheet.Range("G3").Value = "0.0840"
$newSheet.Range("H3").Value = 8

$newSheet.Range("B4").Value = "009975"
$newSheet.Range("C4").Value = "华宝标普美国消费(QDII-LOF)人民币C"
$newSheet.Range("D4").Value = "0.61"
$newSheet.Range("E4").Value = "94.37"
$newSheet.Range("F4").Value = "2.32"
$newSheet.Range("G4").Value = "0.0142"
$newSheet.Range("H4").Value = 8

# Strip the "@" number-format styling back off the text cells (C2:C4 never
# had NumberFormat touched, so they're still on the default/unstyled xf).
$newSheet.Range("C2:C4").Copy()
$newSheet.Range("B2:B4").PasteSpecial(-4122)
$newSheet.Range("D2:D4").PasteSpecial(-4122)
$newSheet.Range("E2:E4").PasteSpecial(-4122)
$newSheet.Range("F2:F4").PasteSpecial(-4122)
$newSheet.Range("G2:G4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: prepend the 2022-Q1 row to the "总计" (total) sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Re-apply the index-column style + a clean (no-style) look for B:D by
# copying formats down from the rows that already have them, so nothing
# picks up the transient formatting the row-insert leaves behind.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.18

# Renumber the index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
